# "Generate Report for Archive"
#
# The localization status report was regenerated: the status for the two
# in-flight items flips from "Ready for handoff" to "In Translation" (this
# is a shared-string value, so every cell that shows that status - on the
# Overview sheet as well as on each per-language sheet - needs to be
# updated so they all resolve back to a single shared-string entry), and
# the two "status" columns shrink to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

# --- Update every cell that currently reads "Ready for handoff" ---------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Shrink the status columns to fit the shorter text -------------------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
